$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet 1 ("ورقة1") - main ledger of transactions (rows keyed on column A)
# -------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("I2").Value  = 27712

$ws1.Range("E3").Value  = 500

$ws1.Range("E5").Value  = 3345
$ws1.Range("F5").Value  = 0
$ws1.Range("I5").Value  = 300

$ws1.Range("E6").Value  = 3289
$ws1.Range("F6").Value  = 3300

$ws1.Range("E7").Value  = 70000

$ws1.Range("F10").Value = 1304

$ws1.Range("B42").Value = 225

$ws1.Range("A45").Value = "علاء المشهراوي"

$ws1.Range("A50").Value = "ابوانس سكيك"
$ws1.Range("E50").Value = 2000

$ws1.Range("E53").Value = 21495
$ws1.Range("H53").Value = 21652

$ws1.Range("E59").Value = 8965

$ws1.Range("F66").Value = 1135

$ws1.Range("E71").Value = 1722
$ws1.Range("F71").Value = 0

$ws1.Range("A74").Value = "عبد الكريم عابدين"
$ws1.Range("E74").Value = 330

$ws1.Range("E76").Value = 4750

$ws1.Range("E77").Value = 133

$ws1.Range("E79").Value = 812

$ws1.Range("F82").Value = 395

$ws1.Range("E83").Value = 0

$ws1.Range("A98").Value = "محمد الزين دليس"
$ws1.Range("H98").Value = 105

$ws1.Range("E104").Value = 0
$ws1.Range("F104").Value = 300

$ws1.Range("E105").Value = 32105

$ws1.Range("A107").Value = "عمار ابوضاهر"
$ws1.Range("E107").Value = 508

$ws1.Range("E113").Value = 0

$ws1.Range("E114").Value = 45349

$ws1.Range("E119").Value = 1357

$ws1.Range("E123").Value = 2880

$ws1.Range("B125").Value = 1120

$ws1.Range("E127").Value = 149276

# -------------------------------------------------------------------------
# Sheet 2 ("ورقة2") - summary ledger (same entries, rows offset by +2)
# -------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("D4").Value   = 27712

$ws2.Range("C5").Value   = -500

$ws2.Range("C7").Value   = -3345
$ws2.Range("D7").Value   = 300

$ws2.Range("C8").Value   = 11

$ws2.Range("C9").Value   = -70000

$ws2.Range("C12").Value  = 1304

$ws2.Range("B44").Value  = -225

$ws2.Range("A47").Value  = "علاء المشهراوي"

$ws2.Range("A52").Value  = "ابوانس سكيك"
$ws2.Range("C52").Value  = -2000

$ws2.Range("C55").Value  = -21495
$ws2.Range("D55").Value  = -21652

$ws2.Range("C61").Value  = -8965

$ws2.Range("C68").Value  = 1135

$ws2.Range("C73").Value  = -1722

$ws2.Range("A76").Value  = "عبد الكريم عابدين"
$ws2.Range("C76").Value  = -330

$ws2.Range("C78").Value  = -4750

$ws2.Range("C79").Value  = -133

$ws2.Range("C81").Value  = -812

$ws2.Range("C84").Value  = 395

$ws2.Range("C85").Value  = 0

$ws2.Range("A100").Value = "محمد الزين دليس"
$ws2.Range("D100").Value = -105

$ws2.Range("C106").Value = 300

$ws2.Range("C107").Value = -32105

$ws2.Range("A109").Value = "عمار ابوضاهر"
$ws2.Range("C109").Value = -508

$ws2.Range("C115").Value = 0

$ws2.Range("C116").Value = -45349

$ws2.Range("C121").Value = -1357

$ws2.Range("C125").Value = -2880

$ws2.Range("B127").Value = -1120

$ws2.Range("C129").Value = -149276

$ws2.Range("B160").Value = 6527
$ws2.Range("B161").Value = 300
$ws2.Range("B162").Value = 198

# -------------------------------------------------------------------------
# Selection / active-sheet bookkeeping (matches the saved view state):
# sheet1 now has J2:J136 selected (no longer scrolled to the old I137 cell),
# while sheet2 keeps its own selection and remains the active/visible tab.
# -------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("J2:J136").Select()

$ws2.Activate()
$ws2.Range("B165").Select()
